# Register New user Field test
# Adds a new "Register_User" worksheet at the end of the workbook with
# Login / Fname / Lname / Pwd / ConfirmPwd columns and a sample data row,
# mirroring the mailto hyperlink pattern used on the other sheets.

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the current last sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Register_User"

# Header row
$ws.Range("A1").Value = "Login"
$ws.Range("B1").Value = "Fname"
$ws.Range("C1").Value = "Lname"
$ws.Range("D1").Value = "Pwd"
$ws.Range("E1").Value = "ConfirmPwd"

# Sample data row
$ws.Range("A2").Value = "UserA01"
$ws.Range("B2").Value = "FnameA"
$ws.Range("C2").Value = "LnameB"
$ws.Range("D2").Value = "Pwd@123"
$ws.Range("E2").Value = "Pwd@123"

# Hyperlink the password cells (same mailto: convention used elsewhere in
# this workbook), then restyle them with the built-in Hyperlink style.
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Pwd@123")
$ws.Range("D2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Pwd@123")
$ws.Range("E2").Style = "Hyperlink"

# Match the saved selection/active cell on the new sheet.
$ws.Range("E2").Select() | Out-Null
